$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Calr", "Tshr", "ECs", 3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 21.33333197335333, 191.99998776018, 0.01566801733960173, 0.01566801733960173),
    @("ECs", "Calr", "Tshr", "FAPs", 3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 203.5419247975038, 1831.877323177534, 0.1494889973608709, 0.149488997360871),
    @("ECs", "Calr", "Tshr", "M2", 3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 16.193603237326, 145.742429135934, 0.011893203397855, 0.011893203397855),
    @("ECs", "Calr", "Tshr", "sCs", 3, 1, 93.32574466666667, 279.977234, 0.2327963689879921, 0.2327963689879922, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 75.90310358924378, 683.127932303194, 0.05574615088966444, 0.05574615088966445),
    @("FAPs", "Calr", "Tshr", "ECs", 3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 37.66929222895667, 339.02363006061, 0.02766577319243982, 0.02766577319243982),
    @("FAPs", "Calr", "Tshr", "FAPs", 3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 359.4037844448492, 3234.634060003643, 0.2639599258865907, 0.2639599258865907),
    @("FAPs", "Calr", "Tshr", "M2", 3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 28.593825537827, 257.344429840443, 0.02100040232307878, 0.02100040232307878),
    @("FAPs", "Calr", "Tshr", "sCs", 3, 1, 164.7897643333334, 494.369293, 0.4110597662007076, 0.4110597662007077, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 134.0257674590792, 1206.231907131713, 0.09843366479859836, 0.09843366479859837),
    @("M2", "Calr", "Tshr", "ECs", 3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 21.90710963662, 197.16398672958, 0.01608942166539422, 0.01608942166539422),
    @("M2", "Calr", "Tshr", "FAPs", 3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 209.0163537396393, 1881.147183656754, 0.1535096279729541, 0.1535096279729542),
    @("M2", "Calr", "Tshr", "M2", 3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 16.629143631906, 149.662292687154, 0.01221308097079579, 0.01221308097079579),
    @("M2", "Calr", "Tshr", "sCs", 3, 1, 95.835818, 287.507454, 0.2390576204784642, 0.2390576204784643, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 77.94458053557933, 701.501224820214, 0.05724548986932009, 0.05724548986932009),
    @("sCs", "Calr", "Tshr", "ECs", 3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 10.72971941411333, 96.56747472702, 0.007880317525615464, 0.007880317525615466),
    @("sCs", "Calr", "Tshr", "FAPs", 3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 102.3725569364251, 921.353012427826, 0.07518633278583668, 0.0751863327858367),
    @("sCs", "Calr", "Tshr", "M2", 3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 8.144663911714, 73.30197520542599, 0.005981753602923219, 0.005981753602923219),
    @("sCs", "Calr", "Tshr", "sCs", 3, 1, 46.93870866666666, 140.816126, 0.117086244332836, 0.117086244332836, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 38.17589325428511, 343.583039288566, 0.02803784041846059, 0.0280378404184606),
)
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
